# Atualização de bases das ligas, do dia: 10-06-2024 às 21:53
#
# Two pairs of rows in the "Portugal Segunda Liga" sheet had their data
# (everything except the running index in column A) swapped between the
# two rows: rows 88/89 and rows 140/141.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B (2) through AD (30) — column A (1) keeps its own row index.
$firstCol = 2
$lastCol = 30

$row1 = 88
$row2 = 89
for ($col = $firstCol; $col -le $lastCol; $col++) {
    $val1 = $ws.Cells.Item($row1, $col).Value2
    $val2 = $ws.Cells.Item($row2, $col).Value2
    $ws.Cells.Item($row1, $col).Value2 = $val2
    $ws.Cells.Item($row2, $col).Value2 = $val1
}

$row1 = 140
$row2 = 141
for ($col = $firstCol; $col -le $lastCol; $col++) {
    $val1 = $ws.Cells.Item($row1, $col).Value2
    $val2 = $ws.Cells.Item($row2, $col).Value2
    $ws.Cells.Item($row1, $col).Value2 = $val2
    $ws.Cells.Item($row2, $col).Value2 = $val1
}
